$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("L7").Value = 24
$ws.Range("L19").Value = 2

$ws.Range("L1").Select()
